# Apply "results with fixed workflow" update:
# For both data sheets (NBR and BAR), the first 4 data rows (old Cutoff
# values 1-4) were dropped, shifting the table up so it now starts at
# Cutoff=5 and runs through Cutoff=19 (15 data rows instead of 19).
# The running index in column A is renumbered 0..14, and the
# Reaction_number counts in column C are refreshed with new values.

$wb = $excel.ActiveWorkbook

$newC = @{
    "NBR" = @(842,835,856,856,843,845,857,856,859,848,843,838,831,830,829)
    "BAR" = @(846,857,822,819,828,815,780,780,767,769,759,761,751,750,749)
}

foreach ($ws in $wb.Worksheets) {
    # Remove the first four data rows (rows 2-5), shifting the rest up.
    $ws.Rows("2:5").Delete()

    $values = $newC[$ws.Name]

    for ($i = 0; $i -lt 15; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 1).Value = $i
        $ws.Cells.Item($row, 3).Value = $values[$i]
    }
}
